$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A2 value from "username" to "user"
$ws.Range("A2").Value = "user"

# Update the selection to A3 (as reflected in the saved view state)
$ws.Range("A3").Select()
